# Insert a new price observation row at row 272 (pushes the existing
# rows 272-359 down to 273-360, carrying their data/formatting with them)
# and populate the new row with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A272").EntireRow.Insert()

$ws.Range("A272").Value = 4
$ws.Range("B272").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C272").Value = "Los Lagos"
$ws.Range("D272").Value = 44985
$ws.Range("E272").Value = 10
$ws.Range("F272").Value = 100112032
$ws.Range("G272").Value = "Zapallo italiano"
$ws.Range("H272").Value = "Sin especificar"
$ws.Range("I272").Value = "Primera"
$ws.Range("J272").Value = 250
$ws.Range("K272").Value = 12000
$ws.Range("L272").Value = 13000
$ws.Range("M272").Value = 12400
$ws.Range("N272").Value = "$/caja 50 unidades"
$ws.Range("O272").Value = "Región de O'Higgins"
$ws.Range("P272").Value = 248
$ws.Range("Q272").Value = 50
$ws.Range("R272").Value = "Hortaliza"
